$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 47, shifting rows 47:146 down to 48:147
$ws.Rows.Item(47).Insert()

# Populate new row 47 with fresh values (copy unchanged columns from what is now row 48 -- same as old row 47)
$ws.Range("A47").Value = 4
$ws.Range("B47").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C47").Value = "Los Lagos"
$ws.Range("D47").Value = 44469
$ws.Range("E47").Value = 10
$ws.Range("F47").Value = 100112021
$ws.Range("G47").Value = "Ají"
$ws.Range("H47").Value = "Inferno"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 90
$ws.Range("K47").Value = 47000
$ws.Range("L47").Value = 47000
$ws.Range("M47").Value = 47000
$ws.Range("N47").Value = "`$/caja 12 kilos"
$ws.Range("O47").Value = "Región de Arica y Parinacota"
$ws.Range("P47").Value = 3917
$ws.Range("Q47").Value = 12
$ws.Range("R47").Value = "Hortaliza"
